$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '26.157.31'
Set-TextValue $ws.Range('E2') '  -4.34%  '
Set-TextValue $ws.Range('D3') '1.656.95'
Set-TextValue $ws.Range('E3') '  -2.97%  '
Set-TextValue $ws.Range('D4') '1.006'
Set-TextValue $ws.Range('E4') '  +0.41%  '
Set-TextValue $ws.Range('D5') '218.02'
Set-TextValue $ws.Range('E5') '  -2.71%  '
Set-TextValue $ws.Range('D6') '0.5174'
Set-TextValue $ws.Range('E6') '  -2.79%  '
Set-TextValue $ws.Range('E7') '  +0.38%  '
Set-TextValue $ws.Range('B8') 'Cardano'
Set-TextValue $ws.Range('C8') 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Range('D8') '0.2572'
Set-TextValue $ws.Range('E8') '  -3.79%  '
Set-TextValue $ws.Range('B9') 'Dogecoin'
Set-TextValue $ws.Range('C9') 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range('D9') '0.06435'
Set-TextValue $ws.Range('E9') '  -2.67%  '
Set-TextValue $ws.Range('E10') '  -4.88%  '
Set-TextValue $ws.Range('D11') '0.07790'
Set-TextValue $ws.Range('D12') '1.664.22'
Set-TextValue $ws.Range('E12') '  -2.67%  '
Set-TextValue $ws.Range('D13') '1.885.48'
Set-TextValue $ws.Range('D14') '4.291'
Set-TextValue $ws.Range('E14') '  -5.67%  '
Set-TextValue $ws.Range('D15') '0.5541'
Set-TextValue $ws.Range('E15') '  -3.93%  '
Set-TextValue $ws.Range('D16') '0.0₅8054'
Set-TextValue $ws.Range('E16') '  -1.54%  '
Set-TextValue $ws.Range('D17') '64.39'
Set-TextValue $ws.Range('E17') '  -4.86%  '
Set-TextValue $ws.Range('D18') '26.211.38'
Set-TextValue $ws.Range('E18') '  -4.08%  '
Set-TextValue $ws.Range('B19') 'Dai'
Set-TextValue $ws.Range('C19') 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D19') '1.006'
Set-TextValue $ws.Range('E19') '  +0.41%  '
Set-TextValue $ws.Range('B20') 'BitcoinCash'
Set-TextValue $ws.Range('C20') 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range('D20') '211.57'
Set-TextValue $ws.Range('E20') '  -2.84%  '
Set-TextValue $ws.Range('D21') '4.397'
Set-TextValue $ws.Range('E21') '  -5.61%  '
Set-TextValue $ws.Range('E22') '  -3.54%  '
Set-TextValue $ws.Range('D23') '5.919'
Set-TextValue $ws.Range('E23') '  -0.67%  '
Set-TextValue $ws.Range('D24') '1.007'
Set-TextValue $ws.Range('E24') '  +0.42%  '
Set-TextValue $ws.Range('D25') '143.92'
Set-TextValue $ws.Range('E25') '  +1.19%  '
Set-TextValue $ws.Range('D26') '1.760'
Set-TextValue $ws.Range('E26') '  +1.82%  '
Set-TextValue $ws.Range('D27') '0.1166'
Set-TextValue $ws.Range('E27') '  -3.78%  '
Set-TextValue $ws.Range('D28') '6.977'
Set-TextValue $ws.Range('D29') '15.77'
Set-TextValue $ws.Range('E29') '  -2.88%  '
Set-TextValue $ws.Range('D30') '0.05279'
Set-TextValue $ws.Range('E30') '  -2.23%  '
Set-TextValue $ws.Range('D31') '1.255'
Set-TextValue $ws.Range('E31') '  -2.71%  '
Set-TextValue $ws.Range('D32') '3.368'
Set-TextValue $ws.Range('E32') '  -3.73%  '
Set-TextValue $ws.Range('D33') '3.230'
Set-TextValue $ws.Range('E33') '  -5.67%  '
Set-TextValue $ws.Range('D34') '1.578'
Set-TextValue $ws.Range('D35') '2.766'
Set-TextValue $ws.Range('E35') '  -3.75%  '
Set-TextValue $ws.Range('D36') '2.365'
Set-TextValue $ws.Range('E36') '  -1.85%  '
Set-TextValue $ws.Range('D37') '0.9254'
Set-TextValue $ws.Range('E37') '  -2.34%  '
Set-TextValue $ws.Range('B38') 'Maker'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D38') '1.170.10'
Set-TextValue $ws.Range('E38') '  +11.92%  '
Set-TextValue $ws.Range('B39') 'ImmutableX'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D39') '0.5720'
Set-TextValue $ws.Range('E39') '  -2.39%  '
Set-TextValue $ws.Range('D40') '0.01594'
Set-TextValue $ws.Range('E40') '  -2.59%  '
Set-TextValue $ws.Range('E41') '  +0.44%  '
Set-TextValue $ws.Range('D42') '0.8370'
Set-TextValue $ws.Range('E42') '  -0.56%  '
Set-TextValue $ws.Range('D43') '5.665'
Set-TextValue $ws.Range('E43') '  -3.09%  '
Set-TextValue $ws.Range('D44') '99.98'
Set-TextValue $ws.Range('E44') '  -0.82%  '
Set-TextValue $ws.Range('D45') '1.795.66'
Set-TextValue $ws.Range('E45') '  -2.97%  '
Set-TextValue $ws.Range('D46') '0.0₈111'
Set-TextValue $ws.Range('E46') '  -5.89%  '
Set-TextValue $ws.Range('D47') '0.4506'
Set-TextValue $ws.Range('E47') '  -0.07%  '
Set-TextValue $ws.Range('D48') '56.00'
Set-TextValue $ws.Range('E48') '  -3.35%  '
Set-TextValue $ws.Range('D49') '1.009'
Set-TextValue $ws.Range('D50') '7.918'
Set-TextValue $ws.Range('E50') '  -2.40%  '
Set-TextValue $ws.Range('D51') '0.05085'
Set-TextValue $ws.Range('E51') '  -2.73%  '
